$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove "Rodrigo" row (old row 6), shifting Matheus/Miguel/Alysson up ---
$ws.Rows("6").Delete()

# --- Insert 9 new columns for meetings 12..20 right before the old "Total" column (N) ---
$ws.Range("N1:V1").EntireColumn.Insert()

# --- Row 3 headers: new meeting columns N..V ---
$ws.Range("N3").Value = "12º Reunião"
$ws.Range("O3").Value = "13º Reunião"
$ws.Range("P3").Value = "14º Reunião"
$ws.Range("Q3").Value = "15º Reunião"
$ws.Range("R3").Value = "16º Reunião"
$ws.Range("S3").Value = "17º Reunião"
$ws.Range("T3").Value = "18º Reunião"
$ws.Range("U3").Value = "19º Reunião"
$ws.Range("V3").Value = "20º Reunião"

# --- Row 4: meeting dates for the new columns + "Presenças" label in the Total column ---
$ws.Range("N4").Value = 44681
$ws.Range("O4").Value = 44688
$ws.Range("P4").Value = 44695
$ws.Range("Q4").Value = 44702
$ws.Range("R4").Value = 44706
$ws.Range("S4").Value = 44707
$ws.Range("T4").Value = 44708
$ws.Range("U4").Value = 44709
$ws.Range("V4").Value = 44710
$ws.Range("W4").Value = "Presenças"

# --- Attendance data: everyone attended every meeting now (was partial before) ---
foreach ($r in 5..8) {
    foreach ($col in @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")) {
        $ws.Range("$col$r").Value = 1
    }
}

# --- Column widths (re-sized for the wider/renumbered table) ---
$ws.Range("B1").EntireColumn.ColumnWidth = 41.140625
$ws.Range("C1").EntireColumn.ColumnWidth = 13.42578125
$ws.Range("D1:K1").EntireColumn.ColumnWidth = 13.85546875
$ws.Range("L1").EntireColumn.ColumnWidth = 14.85546875
$ws.Range("M1").EntireColumn.ColumnWidth = 14.42578125
$ws.Range("N1:U1").EntireColumn.ColumnWidth = 14.85546875
$ws.Range("V1").EntireColumn.ColumnWidth = 15.28515625
$ws.Range("W1").EntireColumn.ColumnWidth = 13.5703125

# --- View: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("E4").Select()

# --- Conditional formatting ranges (grown from column N to column W) ---
$fcs = $ws.Cells.FormatConditions
$fcs.Item(3).ModifyAppliesToRange($ws.Range("C4:W4"))
$fcs.Item(5).Priority = 21
$fcs.Item(6).ModifyAppliesToRange($ws.Range("C3:W3"))
$fcs.Item(6).Priority = 24
$fcs.Item(7).ModifyAppliesToRange($ws.Range("C3:W3"))
$fcs.Item(7).Priority = 26
$fcs.Item(1).ModifyAppliesToRange($ws.Range("C3:W3"))
$fcs.Item(2).ModifyAppliesToRange($ws.Range("C3:W3"))
